# Updated cryptos list on Sun Dec 31 02:43:42 UTC 2023 with GitHub Actions
# Refreshes Price (col D) / Volume(1h) (col E) figures for each coin row,
# and swaps two row-pairs whose relative ranking flipped (Kaspa/WEMIXToken,
# MultiversX/FirstDigitalUSD) by rewriting B/C/D/E in place.
#
# A leading '' (PowerShell-escaped single quote) forces Excel's text/quote-
# prefix interpretation for values that would otherwise auto-parse as a
# number (e.g. "318.25"), matching the original text-formatted cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.334.94'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.287.22'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''318.25'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '''101.31'
$ws.Range("D7").Value = '''0.627'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.602'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").Value = '''39.18'
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").Value = '''0.0902'
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").Value = '''8.25'
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("D13").Value = '''0.106'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '''0.960'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = '''15.17'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").Value = '2.634.78'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '2.296.42'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").Value = '42.288.29'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '''7.40'
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Value = '''13.06'
$ws.Range("E21").Value = '  +33.04%  '
$ws.Range("D22").Value = '''72.76'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '''3.55'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").Value = '''267.28'
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").Value = '''2.21'
$ws.Range("E25").Value = '  -5.06%  '
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").Value = '''10.81'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  +2.56%  '
$ws.Range("D29").Value = '''22.51'
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").Value = '''37.42'
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("D31").Value = '''166.03'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = '''6.11'
$ws.Range("E32").Value = '  +3.76%  '
$ws.Range("D33").Value = '''0.0871'
$ws.Range("E33").Value = '  -2.48%  '
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.55'
$ws.Range("E35").Value = '  -12.70%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.115'
$ws.Range("E36").Value = '  -4.62%  '
$ws.Range("D37").Value = '''4.57'
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("D38").Value = '''0.0358'
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("D39").Value = '''3.65'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("D41").Value = '''1.52'
$ws.Range("E41").Value = '  +2.43%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '''68.52'
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("D45").Value = '''93.06'
$ws.Range("E45").Value = '  -6.95%  '
$ws.Range("D46").Value = '''114.63'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").Value = '''11.93'
$ws.Range("D48").Value = '''78.82'
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("D49").Value = '''8.95'
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").Value = '''5.22'
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").Value = '1.599.87'
$ws.Range("E51").Value = '  +3.25%  '
